$wb = $excel.ActiveWorkbook
$sheet10 = $wb.Worksheets.Item("10")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet10)
$newSheet.Name = "11"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
